$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto data snapshot.
# NumberFormat "@" + Style "Normal" keeps ambiguous numeric-looking strings (e.g. "1.00")
# stored as literal text (matching the pre-existing inline-string cells) without leaving
# a lingering custom style applied to the cell.
$data = @(
    @{Row=2; D="65.517.83"; E="  -1.12%  "},
    @{Row=3; D="3.437.86"; E="  -4.05%  "},
    @{Row=4; E="  -0.17%  "},
    @{Row=5; D="593.80"; E="  -1.91%  "},
    @{Row=6; D="135.50"; E="  -8.04%  "},
    @{Row=7; D="3.436.41"; E="  -3.93%  "},
    @{Row=8; E="  -0.15%  "},
    @{Row=9; E="  +0.49%  "},
    @{Row=10; D="7.54"; E="  -3.64%  "},
    @{Row=11; E="  -9.69%  "},
    @{Row=12; E="  -8.05%  "},
    @{Row=13; D="4.018.00"; E="  -4.16%  "},
    @{Row=14; E="  -11.95%  "},
    @{Row=15; D="26.49"; E="  -9.92%  "},
    @{Row=16; D="3.460.86"; E="  -3.78%  "},
    @{Row=17; D="65.434.70"; E="  -1.42%  "},
    @{Row=19; D="10.01"; E="  -9.18%  "},
    @{Row=20; E="  -8.71%  "},
    @{Row=21; D="13.71"; E="  -7.49%  "},
    @{Row=22; D="392.21"; E="  -7.26%  "},
    @{Row=23; E="  -10.32%  "},
    @{Row=24; D="73.06"; E="  -6.76%  "},
    @{Row=25; E="  +0.05%  "},
    @{Row=26; D="3.579.30"; E="  -4.06%  "},
    @{Row=27; E="  -11.86%  "},
    @{Row=28; E="  +0.07%  "},
    @{Row=29; D="7.32"; E="  -10.49%  "},
    @{Row=30; E="  -9.03%  "},
    @{Row=31; E="  -12.43%  "},
    @{Row=32; D="3.444.25"; E="  -3.85%  "},
    @{Row=33; E="  -0.03%  "},
    @{Row=34; E="  -7.57%  "},
    @{Row=35; D="22.71"; E="  -8.97%  "},
    @{Row=36; D="171.92"; E="  -1.59%  "},
    @{Row=37; D="1.22"; E="  -13.52%  "},
    @{Row=38; D="6.85"; E="  -11.28%  "},
    @{Row=39; E="  -8.34%  "},
    @{Row=40; D="4.79"; E="  -13.57%  "},
    @{Row=41; E="  -8.93%  "},
    @{Row=42; D="0.814"; E="  -7.38%  "},
    @{Row=43; D="43.55"; E="  -5.10%  "},
    @{Row=44; D="1.00"; E="  -0.09%  "},
    @{Row=45; E="  -14.60%  "},
    @{Row=46; E="  -12.45%  "},
    @{Row=47; D="22.96"; E="  -2.46%  "},
    @{Row=48; E="  -1.77%  "},
    @{Row=49; E="  -8.35%  "},
    @{Row=50; E="  -15.17%  "},
    @{Row=51; D="2.190.36"; E="  -8.15%  "}
)

foreach ($item in $data) {
    if ($item.ContainsKey("D")) {
        $dCell = $ws.Range("D" + $item.Row)
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    $eCell = $ws.Range("E" + $item.Row)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.E
    $eCell.Style = "Normal"
}
